$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Soils" sheet: no data changes, only move the cursor/selection to E5.
# ---------------------------------------------------------------------------
$wsSoils = $wb.Worksheets.Item("Soils")
$wsSoils.Range("E5").Select()

# ---------------------------------------------------------------------------
# "Clusters" sheet: update row 2 (IW1/SUB1) with new figures, then drop the
# two extra rows (SUB2, SUB3) that used to follow it.
# ---------------------------------------------------------------------------
$wsClusters = $wb.Worksheets.Item("Clusters")
$wsClusters.Activate()

$wsClusters.Range("D2").Value = 10
$wsClusters.Range("F2").Value = 33.33
$wsClusters.Range("G2").Value = 33.33
$wsClusters.Range("H2").Value = 33.33

$wsClusters.Range("A3:H4").EntireRow.Delete()

$wsClusters.Range("D3").Select()
$excel.ActiveWindow.Zoom = 85

# ---------------------------------------------------------------------------
# "Profiles" sheet: the whole table (header + 5 data rows) is removed.
# ---------------------------------------------------------------------------
$wsProfiles = $wb.Worksheets.Item("Profiles")
$wsProfiles.Activate()

$wsProfiles.Range("A1:D5").EntireRow.Delete()

$wsProfiles.Range("A1:XFD5").Select()
$excel.ActiveWindow.Zoom = 115

# ---------------------------------------------------------------------------
# Leave "Clusters" as the active tab/sheet.
# ---------------------------------------------------------------------------
$wsClusters.Activate()
$wsClusters.Range("D3").Select()
